$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.685.39"
$ws.Range("E2").Value = "  +0.89%  "
$ws.Range("D3").Value = "1.851.13"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "262.87"
$ws.Range("E5").Value = "  -0.91%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5375"
$ws.Range("E7").Value = "  +3.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3203"
$ws.Range("E8").Value = "  -2.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07048"
$ws.Range("E9").Value = "  +3.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.11"
$ws.Range("E10").Value = "  +1.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7795"
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07827"
$ws.Range("E12").Value = "  +0.59%  "
$ws.Range("D13").Value = "1.867.69"
$ws.Range("E13").Value = "  +3.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "89.60"
$ws.Range("E14").Value = "  +1.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.058"
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.17"
$ws.Range("E16").Value = "  +1.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.0000"
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008021"
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").Value = "26.700.12"
$ws.Range("E20").Value = "  +0.92%  "
$ws.Range("D21").Value = "2.089.32"
$ws.Range("E21").Value = "  +1.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.652"
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.055"
$ws.Range("E23").Value = "  +0.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.423"
$ws.Range("E24").Value = "  -1.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.10"
$ws.Range("E25").Value = "  -0.91%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.222"
$ws.Range("E26").Value = "  +1.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.697"
$ws.Range("E27").Value = "  +1.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.17"
$ws.Range("E28").Value = "  +0.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "111.88"
$ws.Range("E29").Value = "  -0.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.300"
$ws.Range("E30").Value = "  +2.80%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.124"
$ws.Range("E31").Value = "  -0.59%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08760"
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04885"
$ws.Range("E33").Value = "  +0.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7411"
$ws.Range("E34").Value = "  +2.76%  "
$ws.Range("E35").Value = "  +0.70%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.882"
$ws.Range("E36").Value = "  +0.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.113"
$ws.Range("E37").Value = "  +0.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.358"
$ws.Range("E38").Value = "  +6.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01756"
$ws.Range("E39").Value = "  -1.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4848"
$ws.Range("E40").Value = "  -0.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9098"
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "109.67"
$ws.Range("E42").Value = "  -1.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.926"
$ws.Range("E43").Value = "  -2.47%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.742"
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4222"
$ws.Range("E46").Value = "  +1.05%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1257"
$ws.Range("E47").Value = "  +1.48%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.133"
$ws.Range("E48").Value = "  +0.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.10"
$ws.Range("E49").Value = "  +0.24%  "
$ws.Range("E50").Value = "  -1.66%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8987"
$ws.Range("E51").Value = "  +0.59%  "
